# Reposition a few shapes on the first slide (resume header layout shift).
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)

$s1.Shapes.Item("AutoShape 7").Top = 2362200 / 12700
$s1.Shapes.Item("TextBox 14").Top  = 2362200 / 12700
$s1.Shapes.Item("TextBox 16").Top  = 2692236 / 12700
$s1.Shapes.Item("TextBox 17").Top  = 4495800 / 12700

# On slides 2 and 3, the purple divider line ("AutoShape 7") and the
# "EXPERIENCIA_RESTO" textbox ("TextBox 16") are moved up (new Top) and
# brought to the very front of the z-order. Duplicating the shape and
# deleting the original reproduces what PowerPoint does internally when
# a shape is brought to front: the shape gets a new id and is appended
# at the end of the shape tree.
foreach ($idx in 2, 3) {
    $s = $p.Slides.Item($idx)

    $line = $s.Shapes.Item("AutoShape 7")
    $newLine = $line.Duplicate()
    $newLine.Left = 3349881 / 12700
    $newLine.Top  = 2362200 / 12700
    $line.Delete()

    $box = $s.Shapes.Item("TextBox 16")
    $newBox = $box.Duplicate()
    # 277.12506103515625 pt is the nearest representable point value that
    # still converts back to the exact target of 3519488 EMU for Left.
    $newBox.Left = 277.12506103515625
    $newBox.Top  = 2692236 / 12700
    $box.Delete()
}
